$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates — values are numeric-looking text stored as
# inline strings in the source sheet, so we prefix with a quote to force
# text entry (avoiding float coercion/rounding), then clear the resulting
# quote-prefix formatting so the cell style matches the original (no style).
$ws.Range("D2").Value = "'264.71"
$ws.Range("D2").ClearFormats()
$ws.Range("D4").Value = "'6.285"
$ws.Range("D4").ClearFormats()
$ws.Range("D5").Value = "'0.06141"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").Value = "'3.600"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").Value = "'6.682"
$ws.Range("D7").ClearFormats()
$ws.Range("D8").Value = "'1.347"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").Value = "'0.8297"
$ws.Range("D9").ClearFormats()
$ws.Range("D11").Value = "'0.1590"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "'0.08242"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").Value = "'0.03421"
$ws.Range("D13").ClearFormats()
$ws.Range("D14").Value = "'0.03101"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").Value = "'0.09253"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").Value = "'3.911"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").Value = "'0.001718"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").Value = "'0.04879"
$ws.Range("D18").ClearFormats()
$ws.Range("D19").Value = "'0.006263"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").Value = "'0.005271"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").Value = "'0.001088"
$ws.Range("D21").ClearFormats()
$ws.Range("D24").Value = "'2.289"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").Value = "'0.3378"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").Value = "'0.1227"
$ws.Range("D26").ClearFormats()
$ws.Range("D27").Value = "'0.0002681"
$ws.Range("D27").ClearFormats()
$ws.Range("D40").Value = "'0.04617"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").Value = "'0.006954"
$ws.Range("D41").ClearFormats()
$ws.Range("D43").Value = "'0.003130"
$ws.Range("D43").ClearFormats()
$ws.Range("D45").Value = "'0.00006152"
$ws.Range("D45").ClearFormats()
$ws.Range("D47").Value = "'0.7783"
$ws.Range("D47").ClearFormats()
$ws.Range("D48").Value = "'0.1951"
$ws.Range("D48").ClearFormats()
$ws.Range("D50").Value = "'0.01241"
$ws.Range("D50").ClearFormats()

# Volume(1h) label (column E) text updates — plain text, no coercion risk.
$ws.Range("E27").Value = '26UpBotsUBXTWorstin24h'
$ws.Range("E44").Value = '43LocalTradersLCT'
